{"js": "// Update the \"three-digit number x one-digit number\" practice table:\n// each cell's multiplication expression is swapped for a new one.\n// Every old value is unique in the document and none of the new values\n// collide with any old value, so a simple ordered search/replace pass\n// is safe.\nconst replacements = [\n  [\"209\u00d78=\", \"119\u00d76=\"],\n  [\"245\u00d72=\", \"173\u00d78=\"],\n  [\"898\u00d73=\", \"408\u00d72=\"],\n  [\"421\u00d77=\", \"501\u00d73=\"],\n  [\"503\u00d74=\", \"715\u00d77=\"],\n  [\"531\u00d72=\", \"241\u00d78=\"],\n  [\"806\u00d79=\", \"614\u00d73=\"],\n  [\"831\u00d72=\", \"480\u00d79=\"],\n  [\"342\u00d74=\", \"655\u00d75=\"],\n  [\"770\u00d73=\", \"621\u00d78=\"],\n  [\"950\u00d74=\", \"980\u00d79=\"],\n  [\"870\u00d74=\", \"692\u00d73=\"],\n  [\"950\u00d79=\", \"579\u00d78=\"],\n  [\"865\u00d77=\", \"553\u00d75=\"],\n  [\"147\u00d75=\", \"302\u00d73=\"],\n  [\"402\u00d75=\", \"872\u00d77=\"],\n  [\"861\u00d72=\", \"831\u00d78=\"],\n  [\"430\u00d78=\", \"207\u00d78=\"],\n  [\"726\u00d72=\", \"896\u00d76=\"],\n  [\"243\u00d73=\", \"521\u00d77=\"],\n  [\"210\u00d75=\", \"439\u00d77=\"],\n  [\"868\u00d76=\", \"377\u00d79=\"],\n  [\"996\u00d78=\", \"766\u00d74=\"],\n  [\"225\u00d78=\", \"622\u00d79=\"],\n  [\"604\u00d77=\", \"272\u00d75=\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\n\nreturn `replaced ${totalReplaced} of ${replacements.length}`;", "ps1": "# Replace the multiplication expressions throughout the document's table\n# with their new values. Each (old, new) pair is unique, so a simple\n# Find/Replace pass for each pair is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{old=\"209\u00d78=\"; new=\"119\u00d76=\"},\n  @{old=\"245\u00d72=\"; new=\"173\u00d78=\"},\n  @{old=\"898\u00d73=\"; new=\"408\u00d72=\"},\n  @{old=\"421\u00d77=\"; new=\"501\u00d73=\"},\n  @{old=\"503\u00d74=\"; new=\"715\u00d77=\"},\n  @{old=\"531\u00d72=\"; new=\"241\u00d78=\"},\n  @{old=\"806\u00d79=\"; new=\"614\u00d73=\"},\n  @{old=\"831\u00d72=\"; new=\"480\u00d79=\"},\n  @{old=\"342\u00d74=\"; new=\"655\u00d75=\"},\n  @{old=\"770\u00d73=\"; new=\"621\u00d78=\"},\n  @{old=\"950\u00d74=\"; new=\"980\u00d79=\"},\n  @{old=\"870\u00d74=\"; new=\"692\u00d73=\"},\n  @{old=\"950\u00d79=\"; new=\"579\u00d78=\"},\n  @{old=\"865\u00d77=\"; new=\"553\u00d75=\"},\n  @{old=\"147\u00d75=\"; new=\"302\u00d73=\"},\n  @{old=\"402\u00d75=\"; new=\"872\u00d77=\"},\n  @{old=\"861\u00d72=\"; new=\"831\u00d78=\"},\n  @{old=\"430\u00d78=\"; new=\"207\u00d78=\"},\n  @{old=\"726\u00d72=\"; new=\"896\u00d76=\"},\n  @{old=\"243\u00d73=\"; new=\"521\u00d77=\"},\n  @{old=\"210\u00d75=\"; new=\"439\u00d77=\"},\n  @{old=\"868\u00d76=\"; new=\"377\u00d79=\"},\n  @{old=\"996\u00d78=\"; new=\"766\u00d74=\"},\n  @{old=\"225\u00d78=\"; new=\"622\u00d79=\"},\n  @{old=\"604\u00d77=\"; new=\"272\u00d75=\"},\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.old\n  $find.Replacement.Text = $pair.new\n  $find.Forward = $true\n  $find.Wrap = 1  # wdFindContinue\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
